$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 346.30768
$ws.Range("I2").Value = 271.22223
$ws.Range("J2").Value = 515.25
$ws.Range("K2").Value = 271.22223
$ws.Range("L2").Value = 515.25
$ws.Range("M2").Value = -158.22223
$ws.Range("N2").Value = -741.25
$ws.Range("H21").Value = 7603.4
$ws.Range("I21").Value = 4504.25
$ws.Range("J21").Value = 20000
$ws.Range("K21").Value = 4504.25
$ws.Range("L21").Value = 20000
$ws.Range("M21").Value = -4036.25
$ws.Range("N21").Value = -20936
$ws.Range("H23").Value = 7603.4
$ws.Range("I23").Value = 4504.25
$ws.Range("J23").Value = 20000
$ws.Range("K23").Value = 4504.25
$ws.Range("L23").Value = 20000
$ws.Range("M23").Value = -4270.25
$ws.Range("N23").Value = -20468
$ws.Range("H38").Value = 328.07693
$ws.Range("J38").Value = 998.5
$ws.Range("L38").Value = 2995.5
$ws.Range("N38").Value = -3739.5
$ws.Range("H58").Value = 2007.1428
$ws.Range("I58").Value = 1416.6666
$ws.Range("J58").Value = 2450
$ws.Range("K58").Value = 4249.9998
$ws.Range("L58").Value = 7350
$ws.Range("M58").Value = -4099.9998
$ws.Range("N58").Value = -7650
$ws.Range("H87").Value = 26400
$ws.Range("J87").Value = 26400
$ws.Range("L87").Value = 26400
$ws.Range("N87").Value = -28896
$ws.Range("H90").Value = 26400
$ws.Range("J90").Value = 26400
$ws.Range("L90").Value = 79200
$ws.Range("N90").Value = -91680
$ws.Range("H137").Value = 2010.2162
$ws.Range("I137").Value = 1975.12
$ws.Range("J137").Value = 2083.3333
$ws.Range("K137").Value = 5925.36
$ws.Range("L137").Value = 6249.999899999999
$ws.Range("M137").Value = -3375.36
$ws.Range("N137").Value = -11349.9999
$ws.Range("H138").Value = 4307.64
$ws.Range("I138").Value = 483.3
$ws.Range("J138").Value = 6857.2
$ws.Range("K138").Value = 1449.9
$ws.Range("L138").Value = 20571.6
$ws.Range("M138").Value = 3690.1
$ws.Range("N138").Value = -30851.6

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7496.02
$ws.Range("I61").Value = 4364.8125
$ws.Range("J61").Value = 13062.611
$ws.Range("K61").Value = 4364.8125
$ws.Range("L61").Value = 13062.611
$ws.Range("M61").Value = -4152.8125
$ws.Range("N61").Value = -13486.611
$ws.Range("H74").Value = 4306.65
$ws.Range("I74").Value = 1938.1212
$ws.Range("J74").Value = 15472.571
$ws.Range("K74").Value = 1938.1212
$ws.Range("L74").Value = 15472.571
$ws.Range("M74").Value = -1064.1212
$ws.Range("N74").Value = -17220.571
$ws.Range("H77").Value = 4306.65
$ws.Range("I77").Value = 1938.1212
$ws.Range("J77").Value = 15472.571
$ws.Range("K77").Value = 9690.606
$ws.Range("L77").Value = 77362.855
$ws.Range("M77").Value = -5322.606
$ws.Range("N77").Value = -86098.855
$ws.Range("H136").Value = 7496.02
$ws.Range("I136").Value = 4364.8125
$ws.Range("J136").Value = 13062.611
$ws.Range("K136").Value = 13094.4375
$ws.Range("L136").Value = 39187.833
$ws.Range("M136").Value = -10544.4375
$ws.Range("N136").Value = -44287.833

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 78000
$ws.Range("J132").Value = 78000
$ws.Range("L132").Value = 78000
$ws.Range("N132").Value = -88120
$ws.Range("H134").Value = 57084.832
$ws.Range("I134").Value = 1262.5385
$ws.Range("J134").Value = 202222.8
$ws.Range("K134").Value = 3787.6155
$ws.Range("L134").Value = 606668.3999999999
$ws.Range("M134").Value = -1252.6155
$ws.Range("N134").Value = -611738.3999999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1708
$ws.Range("I31").Value = 1146.6487
$ws.Range("J31").Value = 5862
$ws.Range("K31").Value = 1146.6487
$ws.Range("L31").Value = 5862
$ws.Range("M31").Value = -851.6487
$ws.Range("N31").Value = -6452
$ws.Range("H34").Value = 1708
$ws.Range("I34").Value = 1146.6487
$ws.Range("J34").Value = 5862
$ws.Range("K34").Value = 1146.6487
$ws.Range("L34").Value = 5862
$ws.Range("M34").Value = -944.6487
$ws.Range("N34").Value = -6266
$ws.Range("H58").Value = 2842936.8
$ws.Range("I58").Value = 4786340.5
$ws.Range("K58").Value = 4786340.5
$ws.Range("M58").Value = -4786137.5
$ws.Range("H60").Value = 11142.6
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 11142.6
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 11142.6
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -12164.6
$ws.Range("H132").Value = 2180.4888
$ws.Range("I132").Value = 2065.0557
$ws.Range("J132").Value = 2642.2222
$ws.Range("K132").Value = 6195.1671
$ws.Range("L132").Value = 7926.6666
$ws.Range("M132").Value = -3665.1671
$ws.Range("N132").Value = -12986.6666
$ws.Range("H134").Value = 2462.5264
$ws.Range("I134").Value = 2120
$ws.Range("K134").Value = 6360
$ws.Range("M134").Value = -3825
$ws.Range("H136").Value = 2842936.8
$ws.Range("I136").Value = 4786340.5
$ws.Range("K136").Value = 14359021.5
$ws.Range("M136").Value = -14356471.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H59").Value = 12273.8
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 12273.8
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 12273.8
$ws.Range("M59").ClearContents()
$ws.Range("N59").Value = -13439.8
$ws.Range("H97").Value = 1946.25
$ws.Range("I97").Value = 2027.5
$ws.Range("J97").Value = 1865
$ws.Range("K97").Value = 2027.5
$ws.Range("L97").Value = 1865
$ws.Range("M97").Value = -1531.5
$ws.Range("N97").Value = -2857
$ws.Range("H126").Value = 2416.4583
$ws.Range("I126").Value = 1899.6875
$ws.Range("K126").Value = 5699.0625
$ws.Range("M126").Value = -3229.0625

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H60").Value = 13011.5
$ws.Range("J60").Value = 19434.5
$ws.Range("L60").Value = 19434.5
$ws.Range("N60").Value = -20452.5
$ws.Range("H136").Value = 5396.8613
$ws.Range("I136").Value = 3460.1765
$ws.Range("J136").Value = 7129.684
$ws.Range("K136").Value = 10380.5295
$ws.Range("L136").Value = 21389.052
$ws.Range("M136").Value = -7830.529500000001
$ws.Range("N136").Value = -26489.052

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 866
$ws.Range("I96").Value = 599
$ws.Range("K96").Value = 599
$ws.Range("M96").Value = 774
$ws.Range("H132").Value = 3709.875
$ws.Range("I132").Value = 3290.5334
$ws.Range("J132").Value = 10000
$ws.Range("K132").Value = 9871.600199999999
$ws.Range("L132").Value = 30000
$ws.Range("M132").Value = -7341.600199999999
$ws.Range("N132").Value = -35060
$ws.Range("H136").Value = 5603.9736
$ws.Range("I136").Value = 1797.6522
$ws.Range("K136").Value = 5392.9566
$ws.Range("M136").Value = -2842.9566
